$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1. Title heading: "Circle Broader View" -> "Circle Language Broader View"
#    Insert a new "Language " run right after the existing "Circle " run.
# -----------------------------------------------------------------------
$titleRange = $d.Content
$found = $titleRange.Find.Execute("Circle ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $titleRange.Collapse(0)   # wdCollapseEnd - move to just after "Circle "
    $titleRange.InsertAfter("Language ")
}

# -----------------------------------------------------------------------
# 2. Remove the leftover "_GoBack" bookmark that wraps the red "*" run
#    right after "opposite" (an editing-session artifact Word leaves
#    behind and subsequently cleans up on the next save).
# -----------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
}
